$d = $word.ActiveDocument

# Locate the "Who did what?" paragraph and insert a brand-new paragraph
# right after it, mirroring its formatting (Times New Roman rPr carried
# on the pPr mark), then fill in the meeting-notes text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Who did what\?") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "Constructing Laplacian pyramid- Richard coding, Dan and Nhung –looking at documents "
